# Update crypto price/volume data (and two name/link swaps) per latest GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''61.163.82'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +7.06%  '

$ws.Range('D3').Value = '''2.640.88'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +9.22%  '

$ws.Range('D4').Value = '''1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.04%  '

$ws.Range('D5').Value = '''513.76'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.45%  '

$ws.Range('D6').Value = '''157.64'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.48%  '

$ws.Range('D7').Value = '''0.612'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.26%  '

$ws.Range('E8').Value = '  -0.24%  '

$ws.Range('D9').Value = '''2.686.97'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +10.17%  '

$ws.Range('D10').Value = '''6.27'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +11.07%  '

$ws.Range('D11').Value = '''0.105'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +5.58%  '

$ws.Range('D12').Value = '''0.350'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +4.07%  '

$ws.Range('E13').Value = '  +0.84%  '

$ws.Range('D14').Value = '''3.112.35'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +9.48%  '

$ws.Range('D15').Value = '''61.082.46'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +6.87%  '

$ws.Range('D16').Value = '''21.89'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +5.31%  '

$ws.Range('E17').Value = '  +5.21%  '

$ws.Range('D18').Value = '''2.683.23'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +10.15%  '

$ws.Range('D19').Value = '''4.81'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.00%  '

$ws.Range('D20').Value = '''351.80'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +7.99%  '

$ws.Range('D21').Value = '''10.53'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +5.32%  '

$ws.Range('D22').Value = '''6.20'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.87%  '

$ws.Range('D23').Value = '''0.997'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.22%  '

$ws.Range('D24').Value = '''60.30'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.53%  '

$ws.Range('D25').Value = '''0.423'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.17%  '

$ws.Range('D26').Value = '''2.778.78'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +10.01%  '

$ws.Range('E27').Value = '  +5.25%  '

$ws.Range('D28').Value = '''0.991'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.88%  '

$ws.Range('D29').Value = '''0.0₃0873'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +11.06%  '

$ws.Range('D30').Value = '''7.57'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.41%  '

$ws.Range('E31').Value = '  -0.12%  '

$ws.Range('D32').Value = '''19.65'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.27%  '

$ws.Range('D33').Value = '''157.19'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.74%  '

$ws.Range('D34').Value = '''1.58'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.08%  '

$ws.Range('E35').Value = '  +8.17%  '

$ws.Range('D36').Value = '''4.05'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +9.09%  '

$ws.Range('D37').Value = '''1.22'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +6.16%  '

$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').Value = '''1.53'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +11.47%  '

$ws.Range('B39').Value = 'Fetch.AI'
$ws.Range('C39').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D39').Value = '''0.874'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.77%  '

$ws.Range('D40').Value = '''310.23'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +16.66%  '

$ws.Range('D41').Value = '''3.79'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +7.20%  '

$ws.Range('D42').Value = '''0.830'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +28.09%  '

$ws.Range('D43').Value = '''35.71'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.70%  '

$ws.Range('E44').Value = '  +9.07%  '

$ws.Range('E45').Value = '  +8.80%  '

$ws.Range('E46').Value = '  +0.16%  '

$ws.Range('D47').Value = '''20.24'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +15.71%  '

$ws.Range('B48').Value = 'FirstDigitalUSD'
$ws.Range('C48').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D48').Value = '''0.995'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.19%  '

$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').Value = '''5.04'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +6.50%  '

$ws.Range('D50').Value = '''0.0237'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.96%  '

$ws.Range('D51').Value = '''2.051.37'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +9.74%  '
